$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set C2:C12 all to the value 3551 (temp solve of RWheel)
$ws.Range("C2:C12").Value = 3551
